# Widen/narrow several slide-title placeholder text boxes.
# Only the shape width (a:ext/@cx) changes; position (a:off) and
# height (a:ext/@cy) stay the same on every slide below.
#
# Target widths were derived from the target EMU values and converted
# to points (EMU / 12700), picking a point value that survives the
# COM Width round-trip (Single-precision) and lands back on the exact
# target EMU after save.

$p = $ppt.ActivePresentation

# Slide 8  - "Assembly - BWA ( + samtools)"  : 6916800 EMU -> 5597100 EMU
$p.Slides.Item(8).Shapes.Item(1).Width = 440.7165832519531

# Slide 9  - "Assembly - Pilon"              : 6916800 EMU -> 3505200 EMU
$p.Slides.Item(9).Shapes.Item(1).Width = 276.0000305175781

# Slide 11 - "Annotation - Trimmomatic + repeatMasker" : 7308000 EMU -> 7406100 EMU
$p.Slides.Item(11).Shapes.Item(1).Width = 583.1575317382812

# Slide 13 - "Functional annotation - eggNOG-mapper"   : 7158300 EMU -> 7295400 EMU
$p.Slides.Item(13).Shapes.Item(1).Width = 574.4409790039062

# Slide 14 - "Counting  reads - STAR and HTSeq"        : 7158300 EMU -> 6409200 EMU
$p.Slides.Item(14).Shapes.Item(1).Width = 504.6614532470703

# Slide 15 - "Differential expression analysis"        : 5394000 EMU -> 5720100 EMU
$p.Slides.Item(15).Shapes.Item(1).Width = 450.401611328125

# Slide 16 - "DE - Between species - PCA"              : 5587500 EMU -> 5868000 EMU
$p.Slides.Item(16).Shapes.Item(1).Width = 462.0472869873047
